# Applies the "added ifoCAST full series evaluation" update:
# The error-metric table (B2:G11) is shifted up by one row (the former
# row N+1 data becomes row N's data), and a brand new evaluation row is
# appended at the bottom (row 11) for the newest matched quarter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for B2:G11 (ME, MAE, MSE, RMSE, SE, N)
$data = @(
    @(0.03429600028756437, 0.7101898576629941, 1.289702394100111, 1.135650647910752, 1.168041881107486, 18),
    @(-0.0941935348587556, 0.7263696581451928, 1.27951414055504, 1.131156107951082, 1.161919457602362, 17),
    @(0.01452822689458977, 0.6449262614402625, 0.6594288688367413, 0.8120522574543718, 0.8385497318331077, 16),
    @(0.1214336542319591, 0.5559322037065559, 0.6436164972988249, 0.8022571266737522, 0.8208469264720462, 15),
    @(0.110076918177882, 0.6331327781440725, 0.8344370595392797, 0.9134752648754534, 0.941050225208, 14),
    @(0.1202220011512185, 0.7253744013351844, 1.112590511420735, 1.054794061142143, 1.090710141917454, 13),
    @(0.1658066239972836, 0.6309408131211106, 0.7997859358034983, 0.8943075174700805, 0.9178794331058292, 12),
    @(0.229040555996801, 0.6860211576689477, 0.8668260477761954, 0.9310349337034542, 0.9464687626294693, 11),
    @(0.2065580353529775, 0.8649197065923389, 1.231080489777346, 1.109540666121502, 1.149112830650219, 10),
    @(0.2295082277971592, 0.7023734988951453, 1.00369417988346, 1.001845387214744, 1.034358580190629, 9)
)

$cols = @("B", "C", "D", "E", "F", "G")
$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range("$($cols[$j])$r").Value = $row[$j]
    }
}
